$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly roll of the Acelga price series: each existing observation (rows 222-333)
# shifts down by one weekly pair (2 rows), a brand-new week's data is written into
# the now-vacated rows 222-223, and the two oldest rows that rolled off the bottom
# (old rows 332-333) are appended as brand-new rows 334-335.
#
# Columns that carry the per-row data: D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg).
# Each tuple: (row, D, J, K, L, M, P)
$rows = @(
    @(222,44606,2000,500,600,550,275),
    @(223,44606,1400,400,450,425,212),
    @(224,44242,3200,450,500,475,238),
    @(225,44242,1600,350,400,375,188),
    @(226,44223,3000,450,500,475,238),
    @(227,44223,1680,350,400,375,188),
    @(228,44515,2200,450,500,475,238),
    @(229,44515,1440,350,400,375,188),
    @(230,44438,3200,450,500,475,238),
    @(231,44438,1540,350,400,375,188),
    @(232,44249,3200,450,500,475,238),
    @(233,44249,1600,350,400,375,188),
    @(234,44566,2400,450,500,475,238),
    @(235,44566,1400,350,400,375,188),
    @(236,44351,3460,500,600,550,275),
    @(237,44351,1680,400,450,425,212),
    @(238,44494,2000,500,600,550,275),
    @(239,44494,1440,400,450,425,212),
    @(240,44237,3400,450,500,475,238),
    @(241,44237,1720,350,400,375,188),
    @(242,44195,2700,400,500,450,225),
    @(243,44195,1600,300,350,325,162),
    @(244,44300,3400,450,500,475,238),
    @(245,44300,1600,350,400,375,188),
    @(246,44578,2400,450,500,475,238),
    @(247,44578,1480,350,400,375,188),
    @(248,44225,3000,450,500,475,238),
    @(249,44225,1600,350,400,375,188),
    @(250,44454,3400,450,500,475,238),
    @(251,44454,1600,350,400,375,188),
    @(252,44281,3400,450,500,475,238),
    @(253,44281,1600,350,400,375,188),
    @(254,44554,2400,450,500,475,238),
    @(255,44554,1460,350,400,375,188),
    @(256,44162,2800,400,500,450,225),
    @(257,44162,1560,300,350,325,162),
    @(258,44568,2600,450,500,475,238),
    @(259,44568,1500,350,400,375,188),
    @(260,44400,3500,450,500,475,238),
    @(261,44400,1600,350,400,375,188),
    @(262,44536,2000,550,600,575,288),
    @(263,44536,1400,450,500,475,238),
    @(264,44193,2800,400,500,450,225),
    @(265,44193,1560,300,350,325,162),
    @(266,44573,2500,450,500,475,238),
    @(267,44573,1460,350,400,375,188),
    @(268,44403,3000,450,500,475,238),
    @(269,44403,1560,350,400,375,188),
    @(270,44319,2800,500,600,550,275),
    @(271,44319,1480,400,450,425,212),
    @(272,44431,3200,450,500,475,238),
    @(273,44431,1560,350,400,375,188),
    @(274,44365,3500,500,600,550,275),
    @(275,44365,1600,400,450,425,212),
    @(276,44426,3500,450,500,475,238),
    @(277,44426,1600,350,400,375,188),
    @(278,44473,2000,450,500,475,238),
    @(279,44473,1400,350,400,375,188),
    @(280,44421,3000,450,500,475,238),
    @(281,44421,1600,350,400,375,188),
    @(282,44445,3200,450,500,475,238),
    @(283,44445,1600,350,400,375,188),
    @(284,44524,2600,550,600,575,288),
    @(285,44524,1400,450,500,475,238),
    @(286,44557,2400,450,500,475,238),
    @(287,44557,1400,350,400,375,188),
    @(288,44489,3100,550,600,575,288),
    @(289,44489,1400,400,450,425,212),
    @(290,44396,2920,450,500,475,238),
    @(291,44396,1500,350,400,375,188),
    @(292,44232,3000,450,500,475,238),
    @(293,44232,1600,350,400,375,188),
    @(294,44279,3400,450,500,475,238),
    @(295,44279,1600,350,400,375,188),
    @(296,44330,3460,500,600,550,275),
    @(297,44330,1640,400,450,425,212),
    @(298,44370,3400,500,600,550,275),
    @(299,44370,1600,400,450,425,212),
    @(300,44487,2000,500,600,550,275),
    @(301,44487,1400,400,450,425,212),
    @(302,44174,2400,400,500,450,225),
    @(303,44174,1500,300,350,325,162),
    @(304,44200,3000,400,500,450,225),
    @(305,44200,1600,300,350,325,162),
    @(306,44221,3000,450,500,475,238),
    @(307,44221,1600,350,400,375,188),
    @(308,44272,3400,450,500,475,238),
    @(309,44272,1600,350,400,375,188),
    @(310,44214,3000,450,500,475,238),
    @(311,44214,1600,350,400,375,188),
    @(312,44312,3000,550,600,575,288),
    @(313,44312,1500,450,500,475,238),
    @(314,44522,2400,550,600,575,288),
    @(315,44522,1400,450,500,475,238),
    @(316,44543,2200,550,600,575,288),
    @(317,44543,1500,450,500,475,238),
    @(318,44167,2800,400,500,450,225),
    @(319,44167,1600,300,350,325,162),
    @(320,44277,2800,450,500,475,238),
    @(321,44277,1400,350,400,375,188),
    @(322,44258,3600,450,500,475,238),
    @(323,44258,1760,350,400,375,188),
    @(324,44349,3400,500,600,550,275),
    @(325,44349,1600,400,450,425,212),
    @(326,44498,2000,500,600,550,275),
    @(327,44498,1360,400,450,425,212),
    @(328,44179,2800,400,500,450,225),
    @(329,44179,1500,300,350,325,162),
    @(330,44335,3460,500,600,550,275),
    @(331,44335,1600,400,450,425,212),
    @(332,44552,2560,450,500,475,238),
    @(333,44552,1480,350,400,375,188),
    @(334,44160,2800,400,500,450,225),
    @(335,44160,1560,300,350,325,162)

)

foreach ($r in $rows) {
    $row = $r[0]
    $ws.Cells.Item($row, 4).Value  = $r[1]   # D - Fecha
    $ws.Cells.Item($row, 10).Value = $r[2]   # J - Volumen
    $ws.Cells.Item($row, 11).Value = $r[3]   # K - Precio minimo
    $ws.Cells.Item($row, 12).Value = $r[4]   # L - Precio maximo
    $ws.Cells.Item($row, 13).Value = $r[5]   # M - Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $r[6]   # P - Precio $/Kg
}

# New rows 334-335 need the rest of their (otherwise-constant) columns populated too,
# copied straight from the row template used throughout this block.
foreach ($row in @(334, 335)) {
    $ws.Cells.Item($row, 1).Value  = 8
    $ws.Cells.Item($row, 2).Value  = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($row, 3).Value  = "Coquimbo"
    $ws.Cells.Item($row, 5).Value  = 4
    $ws.Cells.Item($row, 6).Value  = 100112009
    $ws.Cells.Item($row, 7).Value  = "Acelga"
    $ws.Cells.Item($row, 8).Value  = "Sin especificar"
    $ws.Cells.Item($row, 14).Value = "$/atado 1,5 a 2 kilos"
    $ws.Cells.Item($row, 15).Value = "Provincia del Elquí"
    $ws.Cells.Item($row, 17).Value = 2
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}

$ws.Cells.Item(334, 9).Value = "Primera"
$ws.Cells.Item(335, 9).Value = "Segunda"

# D (Fecha) carries the date-time number format throughout this block; match it
# on the freshly-appended rows the same way it's applied to every row above them.
$ws.Cells.Item(334, 4).NumberFormat = $ws.Cells.Item(332, 4).NumberFormat
$ws.Cells.Item(335, 4).NumberFormat = $ws.Cells.Item(333, 4).NumberFormat
